# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a new
# handback report run:
#   - Status cells flip from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - The "Latest Handback DateTime" timestamps are refreshed
#   - The stale "Error Detail" messages (version mismatch warnings) are
#     cleared now that the handback is in sync
#   - A few columns are widened to better fit the new (longer) status text

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column widths (as persisted in the worksheet XML) are derived from the
# Excel "ColumnWidth" property via: storedWidth = Round((ColumnWidth + 5/6) * 6) / 6
# so to land as close as possible on a desired stored width we back out the
# ColumnWidth to request.
$wideColWidth   = 29.166666666666668   # -> stored width ~= 29.9777050018311
$narrowColWidth = 12.833333333333332   # -> stored width ~= 13.7470531463623

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2:F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth   # column F (de-de)

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2:C3").Value = $newStatus
$wsZhCn.Range("K2:K3").Value = "2016-10-14 07:46:36"
$wsZhCn.Range("P2:P3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = $wideColWidth     # column C (Status)
$wsZhCn.Columns.Item(16).ColumnWidth = $narrowColWidth  # column P (Error Detail)

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2:C3").Value = $newStatus
$wsDeDe.Range("K2:K3").Value = "2016-10-14 07:46:53"
$wsDeDe.Range("P2:P3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth     # column C (Status)
$wsDeDe.Columns.Item(16).ColumnWidth = $narrowColWidth  # column P (Error Detail)
